$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.84"
$ws.Range("E2").Value = "'2.04%"
$ws.Range("D3").Value = "'38.69"
$ws.Range("E3").Value = "'8.46%"
$ws.Range("D4").Value = "'5.087"
$ws.Range("E4").Value = "'0.98%"
$ws.Range("D5").Value = "'0.08165"
$ws.Range("E5").Value = "'3.33%"
$ws.Range("D6").Value = "'1.995"
$ws.Range("E6").Value = "'7.52%"
$ws.Range("D7").Value = "'7.904"
$ws.Range("E7").Value = "'1.53%"
$ws.Range("D8").Value = "'0.9324"
$ws.Range("E8").Value = "'1.44%"
$ws.Range("D9").Value = "'0.1410"
$ws.Range("E9").Value = "'4.79%"
$ws.Range("D10").Value = "'0.1949"
$ws.Range("E10").Value = "'3.57%"
$ws.Range("D11").Value = "'0.09193"
$ws.Range("E11").Value = "'1.64%"
$ws.Range("D12").Value = "'0.03442"
$ws.Range("E12").Value = "'-0.31%"
$ws.Range("D13").Value = "'0.09853"
$ws.Range("E13").Value = "'0.43%"
$ws.Range("D14").Value = "'0.001409"
$ws.Range("E14").Value = "'0.57%"
$ws.Range("D15").Value = "'0.006091"
$ws.Range("E15").Value = "'0.30%"
$ws.Range("D16").Value = "'3.761"
$ws.Range("E16").Value = "'0.93%"
$ws.Range("D17").Value = "'4.183"
$ws.Range("E18").Value = "'4.20%"
$ws.Range("E20").Value = "'0.23%"
$ws.Range("D21").Value = "'4.803"
$ws.Range("E21").Value = "'-7.48%"
$ws.Range("D22").Value = "'0.2454"
$ws.Range("E22").Value = "'11.88%"
$ws.Range("D23").Value = "'0.04468"
$ws.Range("E23").Value = "'1.59%"
$ws.Range("D24").Value = "'0.001240"
$ws.Range("E24").Value = "'0.27%"
$ws.Range("E25").Value = "'-9.27%"
$ws.Range("D27").Value = "'0.0001301"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("D39").Value = "'0.02121"
$ws.Range("E39").Value = "'10.04%"
$ws.Range("D40").Value = "'0.05177"
$ws.Range("E40").Value = "'-2.84%"
$ws.Range("D41").Value = "'0.007459"
$ws.Range("E41").Value = "'-1.86%"
$ws.Range("D42").Value = "'0.01003"
$ws.Range("E42").Value = "'-1.36%"
$ws.Range("E43").Value = "'2.14%"
$ws.Range("D44").Value = "'0.002132"
$ws.Range("E44").Value = "'-0.88%"
$ws.Range("D45").Value = "'0.009686"
$ws.Range("E45").Value = "'-0.61%"
$ws.Range("D46").Value = "'0.00006318"
$ws.Range("E46").Value = "'2.57%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("E48").Value = "'1.94%"
$ws.Range("D49").Value = "'0.001602"
$ws.Range("E49").Value = "'-3.48%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.07%"
